# Applies the cryptos-list refresh described in the commit:
#   "Updated cryptos list on Tue Jul 11 07:52:06 UTC 2023 with GitHub Actions"
#
# Numeric-looking "Price" values (column D) are written with a leading
# apostrophe so Excel keeps them as literal text (matching the original
# inlineStr cells) instead of auto-converting them to numbers, which would
# silently reformat things like "0.000007582" to scientific notation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.560.16"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "1.885.02"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'247.18"
$ws.Range("E5").Value = "  +6.05%  "
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").Value = "'0.4765"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").Value = "'0.2922"
$ws.Range("E8").Value = "  +3.11%  "
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("D10").Value = "'22.07"
$ws.Range("E10").Value = "  +5.56%  "
$ws.Range("B11").Value = "Litecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D11").Value = "'98.06"
$ws.Range("E11").Value = "  +4.97%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.07721"
$ws.Range("E12").Value = "  +0.24%  "
$ws.Range("D13").Value = "'0.7424"
$ws.Range("E13").Value = "  +9.23%  "
$ws.Range("D14").Value = "1.883.25"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "'5.157"
$ws.Range("E15").Value = "  +1.95%  "
$ws.Range("D16").Value = "'275.05"
$ws.Range("E16").Value = "  +3.83%  "
$ws.Range("D17").Value = "30.563.18"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("E18").Value = "  +0.86%  "
$ws.Range("D19").Value = "'0.000007582"
$ws.Range("E19").Value = "  +0.33%  "
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "2.127.45"
$ws.Range("E21").Value = "  +1.12%  "
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "'5.266"
$ws.Range("E23").Value = "  +2.26%  "
$ws.Range("E24").Value = "  +1.77%  "
$ws.Range("E25").Value = "  +0.77%  "
$ws.Range("D26").Value = "'163.62"
$ws.Range("E26").Value = "  -0.89%  "
$ws.Range("E27").Value = "  +2.83%  "
$ws.Range("D28").Value = "'1.953"
$ws.Range("E28").Value = "  +3.73%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "'0.1007"
$ws.Range("E29").Value = "  +2.41%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.370"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "'1.514"
$ws.Range("E31").Value = "  +4.43%  "
$ws.Range("D32").Value = "'4.331"
$ws.Range("E32").Value = "  +2.97%  "
$ws.Range("D33").Value = "'4.127"
$ws.Range("E33").Value = "  +3.76%  "
$ws.Range("D34").Value = "'0.04818"
$ws.Range("E34").Value = "  +3.67%  "
$ws.Range("E35").Value = "  +1.81%  "
$ws.Range("D36").Value = "'0.7022"
$ws.Range("E36").Value = "  +2.64%  "
$ws.Range("D37").Value = "'2.716"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("E38").Value = "  +3.78%  "
$ws.Range("D39").Value = "'2.751"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("D40").Value = "'6.338"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("D41").Value = "'1.994"
$ws.Range("E41").Value = "  +6.09%  "
$ws.Range("D42").Value = "'71.71"
$ws.Range("E42").Value = "  +1.85%  "
$ws.Range("D43").Value = "'0.4234"
$ws.Range("E43").Value = "  +4.70%  "
$ws.Range("D44").Value = "'0.8417"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").Value = "'102.87"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").Value = "'9.348"
$ws.Range("E47").Value = "  +2.63%  "
$ws.Range("D48").Value = "'7.119"
$ws.Range("E48").Value = "  +2.95%  "
$ws.Range("D49").Value = "'35.73"
$ws.Range("E49").Value = "  +4.90%  "
$ws.Range("D50").Value = "'918.78"
$ws.Range("E50").Value = "  -0.64%  "
$ws.Range("E51").Value = "  +4.98%  "
